# Append new entry: 2026-01-10 12:36 JST
# - Update the "取得日時" timestamp on every existing data row
# - Insert a new row (row 6) with a new job entry
# - Keep the two rows that were previously rows 6 and 7 (now rows 7 and 8)
# - Rebuild the hyperlinks on column F so they point at the right rows/URLs

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newTimestamp = "2026-01-10 12:36:10"

# Remove all existing hyperlinks first; they will be recreated after the
# row insert so that they line up with the correct rows.
$ws.Hyperlinks.Delete()

# Insert a new blank row at position 6 - this pushes the old rows 6 and 7
# down to rows 7 and 8.
$ws.Rows.Item(6).Insert()

# --- Row 2 ---------------------------------------------------------------
$ws.Range("A2").Value = $newTimestamp

# --- Row 3 ---------------------------------------------------------------
$ws.Range("A3").Value = $newTimestamp

# --- Row 4 ---------------------------------------------------------------
$ws.Range("A4").Value = $newTimestamp

# --- Row 5 ---------------------------------------------------------------
$ws.Range("A5").Value = $newTimestamp

# --- Row 6 (new) -----------------------------------------------------------
$ws.Range("A6").Value = $newTimestamp
$ws.Range("B6").Value = "チャットボット(授業引き継ぎ支援システム)システム開発のご相談"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5468743"
$ws.Range("G6").Value = 118
$ws.Range("H6").Value = "◆開発,システム開発"

# --- Row 7 (previously row 6) --------------------------------------------
$ws.Range("A7").Value = $newTimestamp
$ws.Range("B7").Value = "製造業DXプロダクト開発のプロダクトマネージャー募集"
$ws.Range("C7").Value = "システム開発"
$ws.Range("D7").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E7").Value = "期限情報なし"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5468432"
$ws.Range("G7").Value = 75
$ws.Range("H7").Value = "◆開発"

# --- Row 8 (previously row 7) --------------------------------------------
$ws.Range("A8").Value = $newTimestamp
$ws.Range("B8").Value = "【緊急対応】インターネットを活用した電話通知システム構築"
$ws.Range("C8").Value = "システム開発"
$ws.Range("D8").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E8").Value = "期限情報なし"
$ws.Range("F8").Value = "https://www.lancers.jp/work/detail/5468565"
$ws.Range("G8").Value = 33

# --- Rebuild the hyperlinks for column F ----------------------------------
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5468493")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5468303")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5217096")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5468677")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5468743")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5468432")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5468565")
